$wb = $excel.ActiveWorkbook

# The "想去人数" (people interested) counts were refreshed for both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same rows).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 2269
    $ws.Range("F4").Value = 382
    $ws.Range("F6").Value = 6440
    $ws.Range("F7").Value = 309
}
